# Update runs/balls/fours/sixes figures for Eoin Morgan (c) vs Kolkata Knight
# Riders across the per-innings rows (2-8). Values are stored as text in the
# sheet (number-looking strings), so we keep them as text by using Excel's
# leading-apostrophe "treat as text" convention when assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @("40", "25", "5", "2")
    3 = @("34", "23", "3", "1")
    4 = @("30", "34", "3", "1")
    5 = @("15", "12", "2", "0")
    6 = @("39", "29", "2", "2")
    7 = @("68", "35", "5", "6")
    8 = @("17", "9",  "2", "1")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("C$row").Value = "'" + $vals[0]
    $ws.Range("D$row").Value = "'" + $vals[1]
    $ws.Range("E$row").Value = "'" + $vals[2]
    $ws.Range("F$row").Value = "'" + $vals[3]
}
